# Restored from revision of admin on 04/15/2020 09:18:19 AM.TEST
# Author: admin. Type: SAVE.
#
# The "From" value of rule R30 (row 10) on the "Rules" sheet changes
# from 18 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
